# Append a new data row (row 2) to the beverage sales sheet.
# Columns: Id, Name, Category, Brand, Quantity Sold, price, Date, Time
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2 ("50"), G2 ("2024-09-15") and H2 ("01:19:17") look like a number / a
# date / a time, so Excel would normally auto-convert them on assignment.
# Mark those cells as Text first so the values are stored verbatim as
# strings, matching the source data feed (which keeps price/date/time as
# plain text rather than typed numbers).
$ws.Range("F2:H2").NumberFormat = "@"

$ws.Range("A2").Value = "707d22be-52ec-4e60-8c20-5f4b21586443"
$ws.Range("B2").Value = "s3Ida"
$ws.Range("C2").Value = "Water"
$ws.Range("D2").Value = "7UP"
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "50"
$ws.Range("G2").Value = "2024-09-15"
$ws.Range("H2").Value = "01:19:17"
